# Update the two data sheets ("展览" and "全部类型") that carry the
# 丽水 / 龙泉 convention-listing rows. Both sheets hold an identical copy of
# the table, so the same edits are applied twice.
#
# Several columns hold numeric- or date-looking text (e.g. "2024-02-07",
# "50") that must stay plain text rather than being auto-converted to a
# real date/number by Excel's input parser. We use Excel's standard
# "text prefix" apostrophe for those so the stored cell keeps its string
# type (mirrors what a user would do typing these values into Excel).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: was 丽水·新年动漫狂欢盛典 -> name cleared, details swapped in
    # from the 龙泉·金沙温泉酒店 event
    $ws.Range("B2").Value = "'2024-02-07"
    $ws.Range("C2").Value = "'"
    $ws.Range("D2").Value = "金沙路26-1号 龙泉金沙温泉酒店"
    $ws.Range("E2").Value = "2024.02.07 10:30-02.07 16:30"
    $ws.Range("F2").Value = 15
    $ws.Range("G2").Value = "'50"
    $ws.Range("I2").Value = "'"
    $ws.Range("J2").Value = "//i2.hdslb.com/bfs/openplatform/202401/rTvQio211704877379770.jpeg"

    # Row 3: was 龙泉·崩X铁X原ONLY -> becomes 丽水·LPJ 现实X次元动漫展
    $ws.Range("B3").Value = "'2024-02-07"
    $ws.Range("C3").Value = "丽水·LPJ 现实X次元动漫展"
    $ws.Range("D3").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E3").Value = "2024.02.07 09:00-02.07 17:00"
    $ws.Range("F3").Value = 263
    $ws.Range("G3").Value = "'45"
    $ws.Range("I3").Value = "'"
    $ws.Range("J3").Value = "//i1.hdslb.com/bfs/openplatform/202311/lP5IkqWn1699431829470.jpeg"

    # Row 4: was 丽水·LPJ 现实X次元动漫展 -> becomes 龙泉·崩X铁X原ONLY
    # (only the name, date format, and link cell change here)
    $ws.Range("B4").Value = "'2024-02-18"
    $ws.Range("C4").Value = "龙泉·崩X铁X原ONLY"
    $ws.Range("I4").Value = "'"
}
